$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 38) with the latest reading.
# Force column A to be stored as plain text so the date-looking string
# ("2025/09/30") is not auto-converted into a date serial number, then
# restore the default "Normal" cell style so no stray formatting is left
# behind (matches the plain, unstyled cells used by the other data rows).
$ws.Range("A38").NumberFormat = "@"
$ws.Range("A38").Value = "2025/09/30"
$ws.Range("A38").Style = "Normal"
$ws.Range("B38").Value = "火"
$ws.Range("C38").Value = 6
$ws.Range("D38").Value = 168
